$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion text in cell A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$oldText = $wsHoja1.Range("A1").Value2
$newText = $oldText -replace [regex]::Escape("1000 Bs = 5.95 = 24403.75 pesos"), "1000 Bs = 5.95 = 24266.62 pesos"
$newText = $newText -replace [regex]::Escape("24403.75 pesos = 5.97 = 968.98 Bs"), "24266.62 pesos = 5.94 = 972.86 Bs"
$wsHoja1.Range("A1").Value = $newText

# --- Sheet "tasas": update numeric rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("O10").Value = 4079
$wsTasas.Range("O12").Value = 163.77
